$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.446.30"
$ws.Range("E2").Value = "  +1.83%  "

# Row 3
$ws.Range("D3").Value = "2.514.77"
$ws.Range("E3").Value = "  +1.42%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'591.73"
$ws.Range("E5").Value = "  +1.13%  "

# Row 6
$ws.Range("D6").Value = "'176.81"
$ws.Range("E6").Value = "  +1.62%  "

# Row 8
$ws.Range("E8").Value = "  +0.31%  "

# Row 9
$ws.Range("D9").Value = "2.514.76"
$ws.Range("E9").Value = "  +1.43%  "

# Row 10
$ws.Range("D10").Value = "'0.144"
$ws.Range("E10").Value = "  +4.00%  "

# Row 11
$ws.Range("E11").Value = "  -1.18%  "

# Row 12
$ws.Range("D12").Value = "'5.00"
$ws.Range("E12").Value = "  +0.80%  "

# Row 13
$ws.Range("D13").Value = "'0.338"
$ws.Range("E13").Value = "  +1.06%  "

# Row 14
$ws.Range("D14").Value = "3.023.44"
$ws.Range("E14").Value = "  +3.07%  "

# Row 15
$ws.Range("D15").Value = "'25.77"
$ws.Range("E15").Value = "  +1.09%  "

# Row 16
$ws.Range("D16").Value = "68.357.93"
$ws.Range("E16").Value = "  +1.84%  "

# Row 17
$ws.Range("D17").Value = "'0.0000170"
$ws.Range("E17").Value = "  -0.21%  "

# Row 18
$ws.Range("D18").Value = "2.525.03"
$ws.Range("E18").Value = "  +1.95%  "

# Row 19
$ws.Range("D19").Value = "'10.98"
$ws.Range("E19").Value = "  +0.11%  "

# Row 20
$ws.Range("D20").Value = "'7.51"
$ws.Range("E20").Value = "  -1.05%  "

# Row 21
$ws.Range("D21").Value = "'351.01"
$ws.Range("E21").Value = "  +0.11%  "

# Row 22
$ws.Range("D22").Value = "'4.17"
$ws.Range("E22").Value = "  +3.45%  "

# Row 23
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'71.25"
$ws.Range("E23").Value = "  +3.21%  "

# Row 24
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.02%  "

# Row 25
$ws.Range("D25").Value = "'4.22"
$ws.Range("E25").Value = "  -0.29%  "

# Row 26
$ws.Range("E26").Value = "  -4.33%  "

# Row 27
$ws.Range("D27").Value = "'9.18"
$ws.Range("E27").Value = "  -0.02%  "

# Row 28
$ws.Range("D28").Value = "2.623.88"
$ws.Range("E28").Value = "  +0.60%  "

# Row 29
$ws.Range("D29").Value = "'0.994"
$ws.Range("E29").Value = "  -0.60%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0900"
$ws.Range("E30").Value = "  -1.12%  "

# Row 31
$ws.Range("D31").Value = "'510.75"
$ws.Range("E31").Value = "  +1.04%  "

# Row 32
$ws.Range("D32").Value = "'7.82"
$ws.Range("E32").Value = "  +0.78%  "

# Row 33
$ws.Range("E33").Value = "  +1.87%  "

# Row 34
$ws.Range("E34").Value = "  +0.64%  "

# Row 35
$ws.Range("E35").Value = "  +0.03%  "

# Row 36
$ws.Range("E36").Value = "  +1.18%  "

# Row 37
$ws.Range("D37").Value = "'163.38"
$ws.Range("E37").Value = "  +1.43%  "

# Row 38
$ws.Range("E38").Value = "  -0.02%  "

# Row 39
$ws.Range("D39").Value = "'18.40"
$ws.Range("E39").Value = "  +1.19%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'1.78"
$ws.Range("E40").Value = "  +5.05%  "

# Row 41
$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D41").Value = "'1.33"
$ws.Range("E41").Value = "  -0.91%  "

# Row 42
$ws.Range("E42").Value = "  -0.01%  "

# Row 43
$ws.Range("D43").Value = "'0.329"
$ws.Range("E43").Value = "  +0.19%  "

# Row 44
$ws.Range("D44").Value = "'4.84"
$ws.Range("E44").Value = "  -0.05%  "

# Row 45
$ws.Range("D45").Value = "'2.42"
$ws.Range("E45").Value = "  +1.27%  "

# Row 46
$ws.Range("D46").Value = "'149.99"
$ws.Range("E46").Value = "  +4.99%  "

# Row 47
$ws.Range("D47").Value = "'3.57"
$ws.Range("E47").Value = "  +2.13%  "

# Row 48
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'0.520"
$ws.Range("E48").Value = "  +0.95%  "

# Row 49
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0259"
$ws.Range("E49").Value = "  -0.62%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0741"
$ws.Range("E50").Value = "  +0.11%  "

# Row 51
$ws.Range("B51").Value = "Optimism"
$ws.Range("C51").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D51").Value = "'1.60"
$ws.Range("E51").Value = "  +0.73%  "
